$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 0.1005870299094374
$ws.Range("I2").Value = 1.642014351456121
$ws.Range("H3").Value = [double]"5.14172098898158e-42"
$ws.Range("I3").Value = -13.58170251854409
$ws.Range("H4").Value = [double]"2.46192811623977e-69"
$ws.Range("I4").Value = 17.60003022547899
$ws.Range("H5").Value = [double]"9.753236417151387e-18"
$ws.Range("I5").Value = 8.576819651182786
$ws.Range("H6").Value = [double]"1.070470095212138e-23"
$ws.Range("I6").Value = 10.03491920089278
$ws.Range("H7").Value = 0.001396946837793676
$ws.Range("I7").Value = 3.195281150425419
$ws.Range("H8").Value = [double]"7.983756059632076e-196"
$ws.Range("I8").Value = 29.85317630367946
$ws.Range("H9").Value = [double]"7.793518645540226e-99"
$ws.Range("I9").Value = -21.10096151109139
$ws.Range("F10").Value = [double]"4.090125119497236e-14"
$ws.Range("H10").Value = [double]"1.626712725324637e-42"
$ws.Range("I10").Value = -13.66572687144085
$ws.Range("H11").Value = [double]"3.762217274443967e-48"
$ws.Range("I11").Value = -14.58005350840661
$ws.Range("H12").Value = [double]"1.283794319129071e-21"
$ws.Range("I12").Value = 9.551057208083092
$ws.Range("H13").Value = [double]"1.349055976310105e-250"
$ws.Range("I13").Value = 33.81122331986262
$ws.Range("H14").Value = [double]"1.070148427334821e-102"
$ws.Range("I14").Value = -21.51739207891101
$ws.Range("H15").Value = [double]"1.740138330085494e-36"
$ws.Range("I15").Value = -12.6152458157533
$ws.Range("H16").Value = [double]"1.0665277478531e-29"
$ws.Range("I16").Value = -11.31818867908856
$ws.Range("H17").Value = [double]"1.196172388658729e-31"
$ws.Range("I17").Value = 11.70538346877259
$ws.Range("H18").Value = 0.0009644551881238153
$ws.Range("I18").Value = -3.300695427568846
$ws.Range("H19").Value = [double]"8.935120371438574e-93"
$ws.Range("I19").Value = 20.43062486268667
$ws.Range("H20").Value = [double]"1.481729373107182e-105"
$ws.Range("I20").Value = 21.82051756860907
$ws.Range("H21").Value = [double]"5.484327550857091e-112"
$ws.Range("I21").Value = 22.4876370959424
$ws.Range("F22").Value = 0.005424722780159434
$ws.Range("H22").Value = 0.0005540473458010006
$ws.Range("I22").Value = 3.453166290299961
$ws.Range("D23").Value = [double]"4.133503021397444e-26"
$ws.Range("F23").Value = [double]"4.313056631039521e-26"
$ws.Range("H23").Value = [double]"1.700408807688419e-97"
$ws.Range("I23").Value = 20.9546876391584
$ws.Range("H24").Value = [double]"9.909403786726777e-57"
$ws.Range("I24").Value = -15.87196205388654
$ws.Range("H25").Value = [double]"7.9312294778736e-19"
$ws.Range("I25").Value = -8.860982532979264
$ws.Range("H26").Value = [double]"1.431882831483472e-28"
$ws.Range("I26").Value = -11.08816917839088
$ws.Range("H27").Value = [double]"4.665513904449939e-20"
$ws.Range("I27").Value = 9.171483989524392
$ws.Range("H28").Value = [double]"2.58109417204881e-159"
$ws.Range("I28").Value = 26.89400467798536
$ws.Range("H29").Value = [double]"4.571296391100679e-31"
$ws.Range("I29").Value = -11.59111739818112
$ws.Range("H30").Value = [double]"1.15132337411586e-12"
$ws.Range("I30").Value = -7.111087232547288
$ws.Range("H31").Value = [double]"6.747721159686214e-12"
$ws.Range("I31").Value = -6.862898418292664
$ws.Range("H32").Value = [double]"1.479426211990881e-26"
$ws.Range("I32").Value = -10.66528635668025
$ws.Range("F33").Value = [double]"4.641747559277947e-07"
$ws.Range("H33").Value = [double]"1.333851400561177e-61"
$ws.Range("I33").Value = -16.56098822601835
$ws.Range("H34").Value = [double]"1.439184655900378e-12"
$ws.Range("I34").Value = 7.080227050948953
$ws.Range("H35").Value = 0.01894535220142964
$ws.Range("I35").Value = 2.346604510927215
$ws.Range("H36").Value = 0.001181611798443139
$ws.Range("I36").Value = 3.243282038831781
$ws.Range("F37").Value = 0.2074807101666617
$ws.Range("H37").Value = 0.1784363945873845
$ws.Range("I37").Value = 1.34558499148247
$ws.Range("H38").Value = [double]"5.158656229513774e-63"
$ws.Range("I38").Value = 16.75554436211863
$ws.Range("H39").Value = [double]"8.663344879867317e-48"
$ws.Range("I39").Value = -14.5230005383045
$ws.Range("H40").Value = [double]"2.663650447286536e-21"
$ws.Range("I40").Value = -9.47515496481191
$ws.Range("H41").Value = [double]"1.894887519049616e-21"
$ws.Range("I41").Value = -9.51064238050604
$ws.Range("H42").Value = [double]"1.933531764003987e-05"
$ws.Range("I42").Value = 4.272431250462653
$ws.Range("H43").Value = [double]"1.533034176882243e-126"
$ws.Range("I43").Value = 23.92888865479044
$ws.Range("H44").Value = [double]"8.913430203429876e-76"
$ws.Range("I44").Value = -18.42100428646772
$ws.Range("H45").Value = [double]"2.134815740252233e-26"
$ws.Range("I45").Value = -10.63114225812502
$ws.Range("H46").Value = [double]"1.518946452831954e-15"
$ws.Range("I46").Value = -7.975394672466073
